# Insert a new weekly price record as row 63 on the active sheet, pushing the
# existing rows 63:81 down to 64:82 (dimension grows from R81 to R82).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 63:81 down by one row to make room for the new record.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly observation.
$ws.Cells.Item(63, 1).Value = 6
$ws.Cells.Item(63, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44711
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 100114007
$ws.Cells.Item(63, 7).Value = "Jengibre"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 380
$ws.Cells.Item(63, 11).Value = 12000
$ws.Cells.Item(63, 12).Value = 13000
$ws.Cells.Item(63, 13).Value = 12605
$ws.Cells.Item(63, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(63, 15).Value = "Perú"
$ws.Cells.Item(63, 16).Value = 970
$ws.Cells.Item(63, 17).Value = 13
$ws.Cells.Item(63, 18).Value = "Hortaliza"
